# The document is one paragraph / one run containing 16 "Label: Value"
# lines, each line separated by a manual line break (<w:br/>). The edit
# re-labels every line as a generic " <field>  <data> " row (mirroring a
# de-identified / templated report) and appends one extra trailing row.
$d = $word.ActiveDocument

$d.Content.Find.Execute('SNOMED: 59000-A-81403 ', $true, $false, $false, $false, $false, $true, 1, $false, ' 欄位  資料 ', 2) | Out-Null
$d.Content.Find.Execute('病史: Bloody stool ', $true, $false, $false, $false, $false, $true, 1, $false, ' ---  --- ', 2) | Out-Null
$d.Content.Find.Execute('診斷: Intestine, large, labeled as "Ascending colon", endoscopic biopsy --- Adenocarcinoma ', $true, $false, $false, $false, $false, $true, 1, $false, ' 診斷資料號  N/A ', 2) | Out-Null
$d.Content.Find.Execute('組織片數: 5 ', $true, $false, $false, $false, $false, $true, 1, $false, ' 病史  N/A ', 2) | Out-Null
$d.Content.Find.Execute('組織尺寸: 0.5 x 0.2 x 0.2 cm ', $true, $false, $false, $false, $false, $true, 1, $false, ' 診斷結果  N/A ', 2) | Out-Null
$d.Content.Find.Execute('組織部位: Colon mucosa ', $true, $false, $false, $false, $false, $true, 1, $false, ' 組織片數  N/A ', 2) | Out-Null
$d.Content.Find.Execute('切片方式: Endoscopic biopsy ', $true, $false, $false, $false, $false, $true, 1, $false, ' 組織尺寸  N/A ', 2) | Out-Null
$d.Content.Find.Execute('處理方式: Fixed in formalin ', $true, $false, $false, $false, $false, $true, 1, $false, ' 組織部位  N/A ', 2) | Out-Null
$d.Content.Find.Execute('組織顏色: Gray white ', $true, $false, $false, $false, $false, $true, 1, $false, ' 切片方式  N/A ', 2) | Out-Null
$d.Content.Find.Execute('組織形狀: Elastic ', $true, $false, $false, $false, $false, $true, 1, $false, ' 處理方式  N/A ', 2) | Out-Null
$d.Content.Find.Execute('顯微鏡檢查: Fragments of necrotic debris and colon mucosa with proliferation and infiltration of irregular hyperchromatic neoplastic glands arranged mainly in complicated tubulo-papillary fashion, a moderately differentiated adenocarcinoma. Remnants suggestive of a pre-existing adenoma are not seen. ', $true, $false, $false, $false, $false, $true, 1, $false, ' 組織顏色  N/A ', 2) | Out-Null
$d.Content.Find.Execute('參考資料: S04-05069 Gall bladder, cholecystectomy --- Acute gangrenous cholecystitis; S01-01737 Skin, nasal bridge, excisional biopsy --- Basal cell carcinoma ', $true, $false, $false, $false, $false, $true, 1, $false, ' 組織形狀  N/A ', 2) | Out-Null
$d.Content.Find.Execute('住院醫師: N/A ', $true, $false, $false, $false, $false, $true, 1, $false, ' 顯微鏡檢查  N/A ', 2) | Out-Null
$d.Content.Find.Execute('病理醫師: Shu-Han Huang, M.D./SWH ', $true, $false, $false, $false, $false, $true, 1, $false, ' 參考資料  N/A ', 2) | Out-Null
$d.Content.Find.Execute('細胞醫檢師: N/A ', $true, $false, $false, $false, $false, $true, 1, $false, ' 住院醫師  N/A ', 2) | Out-Null
$d.Content.Find.Execute('病理專醫字: 病解專醫字第000477號 ', $true, $false, $false, $false, $false, $true, 1, $false, ' 病理醫師  N/A ', 2) | Out-Null

# Append the extra trailing field row ( 病理專醫字  N/A ) that has no corresponding
# line in the original document; insert a line break + new text run at the very
# end of the document content, right before the final paragraph mark.
$endRange = $d.Range($d.Content.End - 2, $d.Content.End - 2)
$endRange.InsertAfter([char]11 + ' 病理專醫字  N/A ')
